# Remove K.Allen from the Chargers Receiving player data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Receiving")

# Find the row containing "K.Allen" in column B and delete the entire row,
# shifting the rows below it up.
$found = $ws.Range("B:B").Find("K.Allen")
if ($found -ne $null) {
    $row = $found.Row
    $ws.Rows.Item($row).Delete()
}

# Make the Receiving sheet the active sheet/tab, matching the saved workbook state.
$ws.Activate()
$ws.Select()
